# DeveloperGuide: update section of UndoRedoStack to UndoRedoCareTaker
# Applies the ModelComponentClassDiagram.pptx edits:
#  - bump the cached datetimeFigureOut field text on the slide master / all
#    custom layouts (3/23/2018 -> 4/5/2018)
#  - enlarge the "Model" component roundRect box
#  - reposition the "<<interface>> ReadOnlyAddressBook" box
#  - add the new "UndoRedoStack" mini class (decision diamond, numbered
#    callout, arrow and rectangle) to the Logic/Model boundary area

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------
# 1) Refresh the cached date field text across the slide master and all
#    custom layouts (the field itself is regenerated by PowerPoint each
#    time the deck is opened/saved on a later date).
# ---------------------------------------------------------------------
$newDate = "4/5/2018"
$master = $p.SlideMaster

for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $shp = $master.Shapes.Item($i)
    if ($shp.Name -like "*Date*") {
        $shp.TextFrame.TextRange.Text = $newDate
    }
}

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $shp = $layout.Shapes.Item($i)
        if ($shp.Name -like "*Date*") {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

# ---------------------------------------------------------------------
# 2) Resize the "Model" rounded-rectangle container so it grows upward
#    and downward to make room for the new UndoRedoStack class.
# ---------------------------------------------------------------------
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.Name -eq "Rectangle 65") {
        $shp.Left = 88.17834645669291
        $shp.Top = 119.37519685039369
        $shp.Width = 589.8216535433071
        $shp.Height = 252
        break
    }
}

# ---------------------------------------------------------------------
# 3) Move the "<<interface>> ReadOnlyAddressBook" box to its new spot.
# ---------------------------------------------------------------------
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.Name -eq "Rectangle 8" -and $shp.TextFrame.TextRange.Text -like "*ReadOnlyAddressBook*") {
        $shp.Left = 214.3637795275591
        $shp.Top = 158.40472440944883
        break
    }
}

# ---------------------------------------------------------------------
# 4) Add the new UndoRedoStack mini class: a small decision diamond, a
#    numbered "1" callout, a purple arrow pointing at the diamond, and
#    the "UndoRedoStack" rectangle itself.
# ---------------------------------------------------------------------

# 4a. Decision diamond - clone an existing "Flowchart: Decision 96" shape
#     so it keeps the same quick-style / fill / line formatting.
$decisionTemplate = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.Name -eq "Flowchart: Decision 96") {
        $decisionTemplate = $shp
        break
    }
}
$decision = $decisionTemplate.Duplicate()
$decision.Name = "Flowchart: Decision 96"
$decision.Left = 164.5937795275591
$decision.Top = 200.02181102362204

# 4b. Numbered callout "1" - clone an existing numbered TextBox.
$textboxTemplate = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.Name -eq "TextBox 53") {
        $textboxTemplate = $shp
        break
    }
}
$callout = $textboxTemplate.Duplicate()
$callout.Name = "TextBox 54"
$callout.Left = 152.4302362204724
$callout.Top = 193.57236220472442

# 4c. Purple arrow connector pointing down into the decision diamond.
$arrow = $s.Shapes.AddConnector(1, 173.42748031496063, 184.13448818897638, 173.46574803149608, 200.32007874015747)
$arrow.Name = "Straight Arrow Connector 55"
$arrow.Line.Visible = $true
$arrow.Line.ForeColor.RGB = 0xA03070
$arrow.Line.Weight = 1.5
$arrow.Line.BeginArrowheadStyle = 1
$arrow.Line.BeginArrowheadWidth = 2
$arrow.Line.BeginArrowheadLength = 2
$arrow.Line.EndArrowheadStyle = 3
$arrow.Line.EndArrowheadWidth = 2
$arrow.Line.EndArrowheadLength = 2
$arrow.Shadow.Visible = $false
$arrow.VerticalFlip = $true

# 4d. "UndoRedoStack" rectangle - clone an existing purple-outlined
#     rectangle class box.
$rectTemplate = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.Name -eq "Rectangle 8" -and $shp.TextFrame.TextRange.Text -eq "AddressBook") {
        $rectTemplate = $shp
        break
    }
}
$undoRedoStack = $rectTemplate.Duplicate()
$undoRedoStack.Name = "Rectangle 8"
$undoRedoStack.Left = 110.6596062992126
$undoRedoStack.Top = 159.0025196850394
$undoRedoStack.Width = 83.7370866141732
$undoRedoStack.Height = 27.303937007874
$undoRedoStack.TextFrame.TextRange.Text = "UndoRedoStack"
